# Apply the edit described by the commit: add a "BO15 - BO10" / "BO10-BO5"
# comparison block (columns H:K) to the "Sheet2" worksheet (the active tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Header row (row 18): two merged-header labels over H:I and J:K
$ws.Range("H18:I18").Merge()
$ws.Range("H18").Value = "BO15 - BO10"

$ws.Range("J18:K18").Merge()
$ws.Range("J18").Value = "BO10-BO5"

# Data rows 19-25: differences between the BO15/BO10/BO5 blocks
for ($r = 19; $r -le 25; $r++) {
    $prev = $r - 7
    $prev2 = $r - 14

    $ws.Range("H$r").Formula = "=C$r-C$prev"
    $ws.Range("I$r").Formula = "=D$r-D$prev"
    $ws.Range("J$r").Formula = "=C$prev-C$prev2"
    $ws.Range("K$r").Formula = "=D$prev-D$prev2"
}

$wb.Save()
